$wb = $excel.ActiveWorkbook

# Remember the current sheet's name before we restructure anything
$oldName = $wb.ActiveSheet.Name

# Add a brand new worksheet (Excel inserts it before the active sheet,
# so the new sheet becomes #1 and the original becomes #2)
$new = $wb.Worksheets.Add()
$old = $wb.Worksheets.Item(2)

# Preserve the page setup from the original sheet
$new.PageSetup.PaperSize = $old.PageSetup.PaperSize
$new.PageSetup.Orientation = $old.PageSetup.Orientation

# Copy the existing Code/Description/Definition table (A1:C8) and paste it
# one column to the right (B1:D8) on the new sheet, values only
$old.Range("A1:C8").Copy()
$new.Range("B1").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# New first column: "Version" header, "1.0" for every data row (as text)
$new.Range("A1").Value = "Version"

$verRange = $new.Range("A2:A8")
$verRange.Formula = "=""1.0"""
$verRange.Copy()
$verRange.PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Drop the old sheet and rename the new one back to the original name
$wb.Worksheets.Item(2).Delete() | Out-Null
$wb.Worksheets.Item(1).Name = $oldName
